$d = $word.ActiveDocument
$count = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($count)
$p.Range.Delete()
